$d = $word.ActiveDocument

# Locate the "KEY ACHIEVEMENTS AND IMPACT" heading paragraph so the edit is
# anchored to that section rather than relying on a fixed paragraph index
# (the same bullet text also appears earlier, under PROFESSIONAL EXPERIENCE).
$headingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd() -eq "KEY ACHIEVEMENTS AND IMPACT") {
        $headingIndex = $i
        break
    }
}

if ($headingIndex -eq -1) {
    throw "Could not find KEY ACHIEVEMENTS AND IMPACT heading"
}

# Section layout (relative to the heading):
#   headingIndex + 0 : "KEY ACHIEVEMENTS AND IMPACT"   (Heading2)
#   headingIndex + 1 : "Impact"                         (Heading3)
#   headingIndex + 2 : bullet 1 - FEC analysis systems
#   headingIndex + 3 : bullet 2 - cloud-based data warehouse
#   headingIndex + 4 : bullet 3 - ETL pipelines
#   headingIndex + 5 : bullet 4 - Trigonometric algorithm (removed)
#   headingIndex + 6 : bullet 5 - redistricting platform  (removed)
#   headingIndex + 7 : bullet 6 - race coding errors

$bullet = [char]0x2022

function Set-ParagraphText($paraIndex, $expectedOldText, $newText) {
    $r = $d.Paragraphs($paraIndex).Range
    $current = $r.Text.Trim()
    if ($current -ne $expectedOldText) {
        throw "Paragraph $paraIndex text mismatch. Expected '$expectedOldText' but found '$current'"
    }
    # Replace everything except the trailing paragraph mark so the
    # paragraph itself (and its formatting) is preserved.
    $textRange = $d.Range($r.Start, $r.End - 1)
    $textRange.Text = $newText
}

Set-ParagraphText ($headingIndex + 2) `
    "$bullet Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion" `
    "$bullet Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%"

Set-ParagraphText ($headingIndex + 3) `
    "$bullet Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy" `
    "$bullet `$4.7M savings enabled nonprofit access"

Set-ParagraphText ($headingIndex + 4) `
    "$bullet Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets" `
    "$bullet Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions"

# Delete the two whole bullet paragraphs that no longer exist, from the
# later index first so the earlier index stays valid.
$expectedRedistricting = "$bullet Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations"
$actualRedistricting = $d.Paragraphs($headingIndex + 6).Range.Text.Trim()
if ($actualRedistricting -ne $expectedRedistricting) {
    throw "Paragraph $($headingIndex + 6) text mismatch. Expected '$expectedRedistricting' but found '$actualRedistricting'"
}
$d.Paragraphs($headingIndex + 6).Range.Delete()

$expectedTrig = "$bullet Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis"
$actualTrig = $d.Paragraphs($headingIndex + 5).Range.Text.Trim()
if ($actualTrig -ne $expectedTrig) {
    throw "Paragraph $($headingIndex + 5) text mismatch. Expected '$expectedTrig' but found '$actualTrig'"
}
$d.Paragraphs($headingIndex + 5).Range.Delete()

Set-ParagraphText ($headingIndex + 5) `
    "$bullet Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%" `
    "$bullet 178% accuracy improvement in racial classification algorithms"
